$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header values for new columns P (16th) and Q (17th) in row 1
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

# Copy the style used by the existing header row cells (e.g. O1) onto the new header cells
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# For every data row (2-25): swap values in columns I, K, M, O and add new columns P, Q
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P = 2
    $ws.Cells.Item($r, 17).Value = 2   # Q = 2
}
